# atualizacao 16 nov 2020
# Append the latest daily readings (through 2020-11-15) to "Diario",
# and the corresponding new monthly summary row to "Mensal".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Mensal": add one new monthly row (row 14)
# ---------------------------------------------------------------
$wsMensal = $wb.Worksheets.Item("Mensal")

$wsMensal.Range("A14").Value = 44150
$wsMensal.Range("B14").Value = 113.48
$wsMensal.Range("C14").Value = 164.93
$wsMensal.Range("D14").Value = -31.19

# Match the date-formatted style already used by column A (copy from the
# row above, which carries style index "2" - numFmt yyyy-mm-dd).
$wsMensal.Range("A13").Copy()
$wsMensal.Range("A14").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet "Diario": add the new daily rows 368-382 (2020-11-01 .. 2020-11-15)
# ---------------------------------------------------------------
$wsDiario = $wb.Worksheets.Item("Diario")

$dailyData = @(
    @(368, 44136, 154.67, 164.93, -6.22),
    @(369, 44137, 128.97, 164.93, -21.8),
    @(370, 44138, 103.91, 164.93, -37),
    @(371, 44139, 91.23, 164.93, -44.69),
    @(372, 44140, 102.51, 164.93, -37.85),
    @(373, 44141, 92.04000000000001, 164.93, -44.19),
    @(374, 44142, 81, 164.93, -50.89),
    @(375, 44143, 73.64, 164.93, -55.35),
    @(376, 44144, 68.89, 164.93, -58.23),
    @(377, 44145, 67.34999999999999, 164.93, -59.16),
    @(378, 44146, 77.09999999999999, 164.93, -53.25),
    @(379, 44147, 134.66, 164.93, -18.35),
    @(380, 44148, 178.67, 164.93, 8.33),
    @(381, 44149, 179.74, 164.93, 8.98),
    @(382, 44150, 167.88, 164.93, 1.79)
)

foreach ($row in $dailyData) {
    $r = $row[0]
    $wsDiario.Range("A$r").Value = $row[1]
    $wsDiario.Range("B$r").Value = $row[2]
    $wsDiario.Range("C$r").Value = $row[3]
    $wsDiario.Range("D$r").Value = $row[4]
}

# Copy the date style (s="2") from the last pre-existing row (367) across
# the whole newly added A368:A382 range in one shot.
$wsDiario.Range("A367").Copy()
$wsDiario.Range("A368:A382").PasteSpecial(-4122)
